$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 117
$ws.Range("I2").Value = 385
$ws.Range("J2").Value = 1586
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 458
$ws.Range("M2").Value = 19
$ws.Range("N2").Value = 257
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 152
$ws.Range("T2").Value = 263
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 2353
$ws.Range("X2").Value = 2324
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 7
